$wb = $excel.ActiveWorkbook

# Sheet1 (展览) updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 674
$ws1.Range("F4").Value = 933
$ws1.Range("F5").Value = 696
$ws1.Range("F6").Value = 829
$ws1.Range("F7").Value = 389
$ws1.Range("F8").Value = 587
$ws1.Range("F9").Value = 123
$ws1.Range("F10").Value = 1191
$ws1.Range("F11").Value = 620
$ws1.Range("F13").Value = 493
$ws1.Range("F15").Value = 330
$ws1.Range("F16").Value = 332
$ws1.Range("F18").Value = 80
$ws1.Range("F19").Value = 546
$ws1.Range("F20").Value = 61
$ws1.Range("F23").Value = 678

# Sheet2 (演出) updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F5").Value = 97
$ws2.Range("F13").Value = 60

# Sheet4 (全部类型) updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F7").Value = 674
$ws4.Range("F8").Value = 933
$ws4.Range("F9").Value = 696
$ws4.Range("F10").Value = 829
$ws4.Range("F11").Value = 389
$ws4.Range("F12").Value = 587
$ws4.Range("F13").Value = 123
$ws4.Range("F14").Value = 1191
$ws4.Range("F15").Value = 620
$ws4.Range("F16").Value = 97
$ws4.Range("F19").Value = 493
$ws4.Range("F22").Value = 330
$ws4.Range("F24").Value = 332
$ws4.Range("F26").Value = 80
$ws4.Range("F29").Value = 546
$ws4.Range("F32").Value = 60
$ws4.Range("F33").Value = 61
$ws4.Range("F36").Value = 678
